$d = $word.ActiveDocument

$key1 = "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"
$key2 = "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"
$key3 = "AIzaSyCP6HmG3IYx3PgBt4CKXVhxkw-n3NRLoYA"

# Locate the paragraph that already holds the first (existing) API key.
$r = $d.Content
$r.Find.Execute($key1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)

# Two blank paragraphs after the first key's paragraph.
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

# New paragraph holding the second (OpenAI) API key.
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter($key2)

$r = $d.Paragraphs.Last.Range

# Two blank paragraphs after the second key's paragraph.
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

# New paragraph holding the third (Google) API key.
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter($key3)
